$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.003208871385164791
$ws.Range("C2").Value = 0.04071648406533734
$ws.Range("D2").Value = 22.3905356188092
$ws.Range("E2").Value = 2195978.878461985
$ws.Range("G2").Value = 2196001.312922959

# Row 3
$ws.Range("B3").Value = 0.6606524410359556
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 22.3905356188092
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("G3").Value = 1157.743882668973

# Row 4
$ws.Range("B4").Value = 0.6606524410359556
$ws.Range("C4").Value = 10.34677158129881
$ws.Range("D4").Value = 3.537761648806719
$ws.Range("E4").Value = 10.19245300693656
$ws.Range("G4").Value = 24.73763867807805
